{"js": "// The hashLookup() description paragraph ends with this unique, distinctive\n// trailing text. We search for it so the insertion point is located\n// robustly (rather than relying on paragraph indices).\nconst anchorText = \"that created the hash.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate the anchor paragraph ending in '\" + anchorText + \"'\");\n}\n\n// Insert the description paragraph right after the anchor paragraph first\n// (it inherits the anchor's normal, non-bold formatting), then insert the\n// bold \"function header\" paragraph before it. This keeps the description\n// paragraph free of any unnecessary explicit bold=\"0\" overrides.\nconst descriptionText =\n  \"This cracker function when called will load a saved table or generate a new one. \" +\n  \"Every uncracked password will then be individualy passed to the table\\u2019s hashLookup() function.\";\nconst descriptionParagraph = anchorParagraph.insertParagraph(descriptionText, Word.InsertLocation.after);\n\nconst headerParagraph = descriptionParagraph.insertParagraph(\"rainbowAttack()\", Word.InsertLocation.before);\nheaderParagraph.font.bold = true;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends with the distinctive phrase \"...that\n# created the hash.\" -- this is the hashLookup() description paragraph\n# that the two new paragraphs (\"rainbowAttack()\" + its description) get\n# inserted after.\n$searchRange = $d.Content\n$searchRange.Find.Execute(\"that created the hash.\") | Out-Null\n$anchorParagraph = $searchRange.Paragraphs(1)\n$anchorRange = $anchorParagraph.Range\n\n# Insert two blank paragraphs right after the anchor paragraph: one for\n# the bold \"rainbowAttack()\" function-header line, one for its\n# description.\n$anchorRange.InsertParagraphAfter()\n$anchorRange.InsertParagraphAfter()\n\n# Re-locate the anchor paragraph (the document structure changed after the\n# inserts above) so we can walk forward to the two new empty paragraphs.\n$searchRange2 = $d.Content\n$searchRange2.Find.Execute(\"that created the hash.\") | Out-Null\n$anchorParagraph2 = $searchRange2.Paragraphs(1)\n$headerParagraph = $anchorParagraph2.Next()\n$descriptionParagraph = $headerParagraph.Next()\n\n# Fill in the bold function-header paragraph. MoveEnd(1, -1) trims the\n# trailing paragraph mark off the range so only the visible text picks up\n# the Bold formatting.\n$headerRange = $headerParagraph.Range\n$headerRange.MoveEnd(1, -1) | Out-Null\n$headerRange.Text = \"rainbowAttack()\"\n$headerRange.Bold = 1\n\n# Fill in the (non-bold) description paragraph.\n$apostrophe = [char]8217\n$descriptionText = \"This cracker function when called will load a saved table or generate a new one. Every uncracked password will then be individualy passed to the table\" + $apostrophe + \"s hashLookup() function.\"\n$descriptionRange = $descriptionParagraph.Range\n$descriptionRange.MoveEnd(1, -1) | Out-Null\n$descriptionRange.Text = $descriptionText\n"}
